$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format target ranges as text first to avoid Excel auto-converting
# numeric-looking strings (e.g. "1.00", "6.00") into numbers, which would
# change their textual representation. Style is reset back to "Normal"
# afterwards so no stray number-format style is left on the cells.
# (Ranges are kept as separate single-area Range objects since the
# multi-area comma syntax is not reliably honoured by NumberFormat here.)
$fmtRange1 = $ws.Range("B32:C33")
$fmtRange2 = $ws.Range("D2:E51")
$fmtRange1.NumberFormat = "@"
$fmtRange2.NumberFormat = "@"

$ws.Range("D2").Value = "93.241.93"
$ws.Range("E2").Value = "  -3.02%  "
$ws.Range("D3").Value = "3.321.71"
$ws.Range("E3").Value = "  -4.22%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "229.76"
$ws.Range("E5").Value = "  -5.55%  "
$ws.Range("D6").Value = "617.68"
$ws.Range("E6").Value = "  -4.31%  "
$ws.Range("D7").Value = "1.36"
$ws.Range("E7").Value = "  -3.78%  "
$ws.Range("D8").Value = "0.383"
$ws.Range("E8").Value = "  -7.03%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "0.918"
$ws.Range("E10").Value = "  -8.39%  "
$ws.Range("D11").Value = "3.317.13"
$ws.Range("E11").Value = "  -4.36%  "
$ws.Range("D12").Value = "41.47"
$ws.Range("E12").Value = "  -4.46%  "
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").Value = "92.993.63"
$ws.Range("E15").Value = "  -2.99%  "
$ws.Range("D16").Value = "3.951.26"
$ws.Range("E16").Value = "  -3.67%  "
$ws.Range("E17").Value = "  -4.84%  "
$ws.Range("D18").Value = "7.94"
$ws.Range("E18").Value = "  -8.05%  "
$ws.Range("D19").Value = "3.324.82"
$ws.Range("E19").Value = "  -3.73%  "
$ws.Range("D20").Value = "17.09"
$ws.Range("E20").Value = "  -7.24%  "
$ws.Range("D21").Value = "10.89"
$ws.Range("E21").Value = "  -9.02%  "
$ws.Range("D22").Value = "3.41"
$ws.Range("E22").Value = "  +3.89%  "
$ws.Range("D23").Value = "489.38"
$ws.Range("E23").Value = "  -4.92%  "
$ws.Range("D24").Value = "0.448"
$ws.Range("E24").Value = "  -9.82%  "
$ws.Range("D25").Value = "0.0000180"
$ws.Range("E25").Value = "  -7.26%  "
$ws.Range("D26").Value = "6.01"
$ws.Range("E26").Value = "  -8.49%  "
$ws.Range("D27").Value = "89.40"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("D28").Value = "3.501.89"
$ws.Range("E28").Value = "  -3.66%  "
$ws.Range("D29").Value = "11.52"
$ws.Range("E29").Value = "  -7.15%  "
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").Value = "10.97"
$ws.Range("E31").Value = "  -8.58%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.135"
$ws.Range("E32").Value = "  -3.47%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "2.63"
$ws.Range("E33").Value = "  -5.18%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "0.171"
$ws.Range("E35").Value = "  -7.84%  "
$ws.Range("D36").Value = "28.19"
$ws.Range("E36").Value = "  -9.67%  "
$ws.Range("D37").Value = "0.524"
$ws.Range("E37").Value = "  -10.51%  "
$ws.Range("D38").Value = "523.34"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "7.28"
$ws.Range("E40").Value = "  -7.58%  "
$ws.Range("D41").Value = "0.146"
$ws.Range("E41").Value = "  -4.52%  "
$ws.Range("D42").Value = "1.35"
$ws.Range("E42").Value = "  -8.35%  "
$ws.Range("D43").Value = "0.864"
$ws.Range("E43").Value = "  -6.98%  "
$ws.Range("D44").Value = "24.04"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("E45").Value = "  -4.17%  "
$ws.Range("D46").Value = "3.52"
$ws.Range("E46").Value = "  -3.33%  "
$ws.Range("D47").Value = "0.0404"
$ws.Range("E47").Value = "  -4.65%  "
$ws.Range("D48").Value = "5.29"
$ws.Range("E48").Value = "  -5.54%  "
$ws.Range("D49").Value = "52.67"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("E50").Value = "  -4.66%  "
$ws.Range("D51").Value = "7.82"
$ws.Range("E51").Value = "  -5.64%  "

$fmtRange1.Style = "Normal"
$fmtRange2.Style = "Normal"
